$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new day's row of expense data (row 6), matching the style
# (date number format) used by the existing date column cells.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$rowData = New-Object 'object[,]' 1,13
$rowData[0,0]  = 43794   # A6 - Date
$rowData[0,1]  = 0       # B6 - Bread
$rowData[0,2]  = 0       # C6 - Veggies
$rowData[0,3]  = 0       # D6 - Fruits
$rowData[0,4]  = 0       # E6 - Meat
$rowData[0,5]  = 12.5    # F6 - Buzz
$rowData[0,6]  = 3.5     # G6 - Dairy
$rowData[0,7]  = 0       # H6 - Other groceries
$rowData[0,8]  = 0       # I6 - Gas
$rowData[0,9]  = 776     # J6 - Bills
$rowData[0,10] = 1       # K6 - Pisi work day
$rowData[0,11] = 3       # L6 - Beers
$rowData[0,12] = 3       # M6 - Wine glasses

$ws.Range("A6:M6").Value = $rowData

# Move the active selection past the newly entered data, as Excel would
# leave it after the user tabbed through the row.
$ws.Range("N6").Select()
